# IKAnimator position lerp finished
# Add 6 new log rows (36-41) to the Workhours sheet, continuing the
# existing "Task 2" table that ends at row 34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: From, To, Task description
$newRows = @(
    @{ Row = 36; From = 0.47222222222222227; To = 0.4861111111111111;  Text = "Solved IK bug." },
    @{ Row = 37; From = 0.4861111111111111;  To = 0.52083333333333337; Text = "Setting up nicer IK positions." },
    @{ Row = 38; From = 0.52083333333333337; To = 0.53125;             Text = "Lerping the rotation sometimes results in the Ik taking the longer rotation. Searching for a fix." },
    @{ Row = 39; From = 0.5493055555555556;  To = 0.60763888888888895; Text = "Fixing skipping if the goal is changed mid transition." },
    @{ Row = 40; From = 0.60763888888888895; To = 0.65972222222222221; Text = "Make it possible to Deatach the ik from the goal, so it will not teleport after the goal, just treat it as a new goal." },
    @{ Row = 41; From = 0.73958333333333337; To = 0.78819444444444453; Text = "Ik working correctly for every limb, except the rotation lerp." }
)

# Rows that need a taller, two-line height like existing wrapped rows.
$tallRows = @(38, 40)

foreach ($item in $newRows) {
    $r = $item.Row

    # Clone formatting (number format, alignment, borders, etc.) from the
    # last existing data row (33) of the same table, so the new rows look
    # exactly like the ones above them.
    $ws.Range("B33:D33").Copy($ws.Range("B$r`:D$r"))

    $ws.Cells.Item($r, 2).Value = $item.From
    $ws.Cells.Item($r, 3).Value = $item.To
    $ws.Cells.Item($r, 4).Value = $item.Text

    if ($tallRows -contains $r) {
        $ws.Rows($r).RowHeight = 30
    }
}

# Update the view state to match where the user ended up after typing.
$ws.Range("D42").Select()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
